# Auto-generated Excel COM-interop script updating the cryptos list
# (prices + hourly volume %) per the commit "Updated cryptos list on
# Tue Jul 18 19:38:43 UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds free-form text (e.g. "29.825.13",
# "0.06830") that must stay text, not be reinterpreted as a number.
# Force the text number format before writing so values like
# "1.001"/"0.06840" keep their exact literal digits/trailing zeros.
$priceCells = @(
    "D2",
    "D3",
    "D4",
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D20",
    "D21",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.839.22"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "1.899.51"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "0.7633"
$ws.Range("E5").Value = "  +4.14%  "

$ws.Range("D6").Value = "240.36"
$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "0.3062"
$ws.Range("E8").Value = "  -1.38%  "

$ws.Range("D9").Value = "25.47"
$ws.Range("E9").Value = "  -2.82%  "

$ws.Range("D10").Value = "0.06840"
$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("D11").Value = "0.07975"
$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").Value = "1.913.01"
$ws.Range("E12").Value = "  +1.10%  "

$ws.Range("D13").Value = "0.7414"
$ws.Range("E13").Value = "  -3.54%  "

$ws.Range("D14").Value = "5.154"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").Value = "90.97"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").Value = "29.839.21"
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").Value = "13.89"
$ws.Range("E17").Value = "  -1.97%  "

$ws.Range("D18").Value = "5.925"
$ws.Range("E18").Value = "  +2.87%  "

$ws.Range("D19").Value = "242.74"
$ws.Range("E19").Value = "  +1.35%  "

$ws.Range("D20").Value = "0.000007702"
$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "6.937"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").Value = "166.62"
$ws.Range("E24").Value = "  +1.35%  "

$ws.Range("D25").Value = "9.215"
$ws.Range("E25").Value = "  -0.64%  "

$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").Value = "0.1299"
$ws.Range("E27").Value = "  +2.23%  "

$ws.Range("D28").Value = "2.027"
$ws.Range("E28").Value = "  +0.72%  "

$ws.Range("D29").Value = "1.407"
$ws.Range("E29").Value = "  +3.56%  "

$ws.Range("D30").Value = "1.517"
$ws.Range("E30").Value = "  -0.82%  "

$ws.Range("D31").Value = "4.254"
$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("D32").Value = "4.087"
$ws.Range("E32").Value = "  +0.15%  "

$ws.Range("D33").Value = "0.05245"
$ws.Range("E33").Value = "  +3.17%  "

$ws.Range("D34").Value = "1.253"
$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("D35").Value = "0.7256"
$ws.Range("E35").Value = "  -1.39%  "

$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("D37").Value = "0.01919"
$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").Value = "2.778"
$ws.Range("E38").Value = "  +0.19%  "

$ws.Range("D39").Value = "6.159"
$ws.Range("E39").Value = "  -2.18%  "

$ws.Range("D40").Value = "0.4415"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("D41").Value = "71.95"
$ws.Range("E41").Value = "  -3.45%  "

$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("D43").Value = "0.8317"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").Value = "1.880"
$ws.Range("E44").Value = "  -2.56%  "

$ws.Range("D45").Value = "7.618"
$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "99.91"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.765"
$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("D48").Value = "2.050.92"
$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("D49").Value = "36.05"
$ws.Range("E49").Value = "  -2.49%  "

$ws.Range("D50").Value = "1.477"
$ws.Range("E50").Value = "  +1.46%  "

$ws.Range("D51").Value = "0.05935"
$ws.Range("E51").Value = "  -0.12%  "
